$wb = $excel.ActiveWorkbook

# Update the "想去人数" (want-to-go count) values on the "展览" sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 104
$ws1.Range("F3").Value = 940

# Update the same values on the "全部类型" sheet (mirrors the 展览 data)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 104
$ws4.Range("F3").Value = 940
